$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for columns B (TB), C (d2S), D (K), E (IP) for rows 2-24
$beValues = @(
    @(0.7287194209349384,1.65323645889881,0.1529057820181812,0.4998867070740569),
    @(0.3464964993005633,1.65323645889881,0.7127328510149897,0.4998867070740569),
    @(3.182878228561681,1.65323645889881,0.7127328510149897,0.4998867070740569),
    @(1.505614041169197,1.65323645889881,0.7127328510149897,0.4998867070740569),
    @(1.505614041169197,1.65323645889881,0.7127328510149897,6.48142807727062),
    @(3.182878228561681,1.65323645889881,0.1529057820181812,0.4998867070740569),
    @(1.505614041169197,1.65323645889881,0.1529057820181812,6.48142807727062),
    @(0.006876353814593728,0.05231270169004087,0.1529057820181812,0.4998867070740569),
    @(0.7287194209349384,1.65323645889881,0.7127328510149897,0.4998867070740569),
    @(0.7287194209349384,0.3375848360084654,0.7127328510149897,0.4998867070740569),
    @(3.182878228561681,1.65323645889881,0.1529057820181812,0.4998867070740569),
    @(1.505614041169197,1.65323645889881,0.7127328510149897,0.4998867070740569),
    @(0.7287194209349384,0.3375848360084654,16.98373111632243,0.4998867070740569),
    @(3.182878228561681,1.65323645889881,3.082599426703578,0.4998867070740569),
    @(3.182878228561681,1.65323645889881,16.98373111632243,0.4998867070740569),
    @(0.06328177979961902,0.05231270169004087,0.1529057820181812,0.4998867070740569),
    @(3.182878228561681,1.65323645889881,0.7127328510149897,0.4998867070740569),
    @(3.182878228561681,1.65323645889881,0.7127328510149897,0.4998867070740569),
    @(1.505614041169197,0.3375848360084654,0.1529057820181812,0.4998867070740569),
    @(1.505614041169197,1.65323645889881,0.1529057820181812,0.4998867070740569),
    @(3.182878228561681,1.65323645889881,3.082599426703578,0.4998867070740569),
    @(3.182878228561681,1.65323645889881,0.1529057820181812,0.4998867070740569),
    @(1.505614041169197,1.65323645889881,16.98373111632243,0.4998867070740569)
)

$gValues = @(
    3.034748368925986,
    3.21235251628842,
    6.048734245549538,
    4.371470058157054,
    10.35301142835362,
    5.488907176552729,
    9.793184359356808,
    0.7119815445968727,
    3.594575437922795,
    2.27892381503245,
    5.488907176552729,
    4.371470058157054,
    18.54992208033989,
    8.418600821238126,
    22.31973251085698,
    0.768386970581898,
    6.048734245549538,
    6.048734245549538,
    2.495991366269901,
    3.811642989160245,
    8.418600821238126,
    5.488907176552729,
    20.64246832346449
)

for ($i = 0; $i -lt $beValues.Count; $i++) {
    $row = $i + 2
    $vals = $beValues[$i]
    $ws.Range("B$row").Value = $vals[0]
    $ws.Range("C$row").Value = $vals[1]
    $ws.Range("D$row").Value = $vals[2]
    $ws.Range("E$row").Value = $vals[3]
    $ws.Range("G$row").Value = $gValues[$i]
}
